# Updates the cryptos list table (columns B-E, rows 2-51) on the active worksheet
# to reflect the latest scraped coinranking.com data, matching the
# "Updated cryptos list ... with GitHub Actions" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row, column letter, new value, whether the value must be
# forced to Text format so Excel does not auto-convert a numeric-looking
# string (e.g. "0.9998") into a real number.
$updates = @(
    @(2, "D", "30.456.43", $false),
    @(2, "E", "  +0.56%  ", $false),
    @(3, "D", "1.870.58", $false),
    @(3, "E", "  +0.16%  ", $false),
    @(4, "D", "0.9998", $true),
    @(4, "E", "  +0.03%  ", $false),
    @(5, "D", "247.01", $true),
    @(5, "E", "  +1.71%  ", $false),
    @(6, "D", "0.9998", $true),
    @(6, "E", "  +0.02%  ", $false),
    @(7, "D", "0.4737", $true),
    @(7, "E", "  +0.28%  ", $false),
    @(8, "D", "0.2920", $true),
    @(8, "E", "  +1.71%  ", $false),
    @(9, "D", "0.06496", $true),
    @(9, "E", "  +0.31%  ", $false),
    @(10, "D", "22.22", $true),
    @(10, "E", "  +6.54%  ", $false),
    @(11, "D", "0.07720", $true),
    @(11, "E", "  -0.03%  ", $false),
    @(12, "D", "97.74", $true),
    @(12, "E", "  +2.64%  ", $false),
    @(13, "D", "0.7434", $true),
    @(13, "E", "  +5.08%  ", $false),
    @(14, "D", "1.872.58", $false),
    @(14, "E", "  +0.25%  ", $false),
    @(15, "D", "5.151", $true),
    @(15, "E", "  +1.18%  ", $false),
    @(16, "D", "273.77", $true),
    @(16, "E", "  +1.57%  ", $false),
    @(17, "D", "30.442.10", $false),
    @(17, "E", "  +0.63%  ", $false),
    @(18, "D", "13.44", $true),
    @(18, "E", "  +0.61%  ", $false),
    @(19, "D", "0.000007545", $true),
    @(19, "E", "  +0.13%  ", $false),
    @(20, "D", "0.9995", $true),
    @(20, "E", "  -0.07%  ", $false),
    @(21, "D", "2.116.82", $false),
    @(21, "E", "  +0.43%  ", $false),
    @(22, "D", "0.9996", $true),
    @(22, "E", "  +0.02%  ", $false),
    @(23, "D", "5.242", $true),
    @(23, "E", "  +0.64%  ", $false),
    @(24, "D", "6.182", $true),
    @(24, "E", "  +0.85%  ", $false),
    @(25, "D", "9.295", $true),
    @(25, "E", "  -0.48%  ", $false),
    @(26, "D", "163.39", $true),
    @(26, "E", "  -1.29%  ", $false),
    @(27, "D", "18.80", $true),
    @(27, "E", "  -0.28%  ", $false),
    @(28, "D", "1.927", $true),
    @(28, "E", "  +0.66%  ", $false),
    @(29, "D", "0.1004", $true),
    @(29, "E", "  +1.80%  ", $false),
    @(30, "D", "1.365", $true),
    @(30, "E", "  -1.69%  ", $false),
    @(31, "E", "  -0.11%  ", $false),
    @(32, "D", "4.285", $true),
    @(32, "E", "  +0.83%  ", $false),
    @(33, "D", "4.130", $true),
    @(33, "E", "  +2.85%  ", $false),
    @(34, "D", "0.04834", $true),
    @(34, "E", "  +2.17%  ", $false),
    @(35, "D", "1.121", $true),
    @(35, "E", "  +0.13%  ", $false),
    @(36, "D", "0.6963", $true),
    @(36, "E", "  +0.83%  ", $false),
    @(37, "B", "Frax", $false),
    @(37, "C", "https://coinranking.com/coin/KfWtaeV1W+frax-frax", $false),
    @(37, "D", "0.9993", $true),
    @(37, "E", "  +0.04%  ", $false),
    @(38, "B", "HuobiToken", $false),
    @(38, "C", "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht", $false),
    @(38, "D", "2.715", $true),
    @(38, "E", "  +0.43%  ", $false),
    @(39, "B", "VeChain", $false),
    @(39, "C", "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet", $false),
    @(39, "D", "0.01854", $true),
    @(39, "E", "  +0.63%  ", $false),
    @(40, "B", "MXToken", $false),
    @(40, "C", "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx", $false),
    @(40, "D", "2.740", $true),
    @(40, "E", "  +0.32%  ", $false),
    @(41, "B", "FraxShare", $false),
    @(41, "C", "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs", $false),
    @(41, "D", "6.310", $true),
    @(41, "E", "  -0.35%  ", $false),
    @(42, "B", "Aave", $false),
    @(42, "C", "https://coinranking.com/coin/ixgUfzmLR+aave-aave", $false),
    @(42, "D", "72.93", $true),
    @(42, "E", "  +3.52%  ", $false),
    @(43, "B", "RenderToken", $false),
    @(43, "C", "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr", $false),
    @(43, "D", "1.979", $true),
    @(43, "E", "  +4.31%  ", $false),
    @(44, "B", "TheSandbox", $false),
    @(44, "C", "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand", $false),
    @(44, "D", "0.4198", $true),
    @(44, "E", "  +3.00%  ", $false),
    @(45, "B", "PaxDollar", $false),
    @(45, "C", "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp", $false),
    @(45, "D", "0.9995", $true),
    @(45, "E", "  +0.01%  ", $false),
    @(46, "B", "TrustWalletToken", $false),
    @(46, "C", "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt", $false),
    @(46, "D", "0.8385", $true),
    @(46, "E", "  -0.31%  ", $false),
    @(47, "B", "Quant", $false),
    @(47, "C", "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt", $false),
    @(47, "D", "101.82", $true),
    @(47, "E", "  -0.17%  ", $false),
    @(48, "B", "EnergySwap", $false),
    @(48, "C", "https://coinranking.com/coin/SbWqqTui-+energyswap-ens", $false),
    @(48, "D", "9.330", $true),
    @(48, "E", "  +1.29%  ", $false),
    @(49, "B", "Elrond", $false),
    @(49, "C", "https://coinranking.com/coin/omwkOTglq+elrond-egld", $false),
    @(49, "D", "35.48", $true),
    @(49, "E", "  +2.12%  ", $false),
    @(50, "B", "Aptos", $false),
    @(50, "C", "https://coinranking.com/coin/HGYj5JCv5+aptos-apt", $false),
    @(50, "D", "6.996", $true),
    @(50, "E", "  -1.29%  ", $false),
    @(51, "B", "Maker", $false),
    @(51, "C", "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr", $false),
    @(51, "D", "922.47", $true),
    @(51, "E", "  -1.36%  ", $false),
)

foreach ($u in $updates) {
    $r = $u[0]
    $col = $u[1]
    $val = $u[2]
    $forceText = $u[3]
    $cell = $ws.Range("$col$r")
    if ($forceText) {
        $cell.NumberFormat = "@"
    }
    $cell.Value = $val
}
